$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Build the two new border styles exactly once, on sheet 1's merged header ---
# C1 -> thin top+bottom only (no left/right)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Weight = 2
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

# D1 -> thin top+bottom+right (no left)
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Weight = 2
$d1.Borders.Item(7).LineStyle = -4142

# --- Re-use those exact formats everywhere else they are needed,
#     via copy/paste-special so no further style objects get minted ---
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)

$c1.Copy()
$ws2.Range("F1").PasteSpecial(-4122)
$d1.Copy()
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Anonymize "fedcore" -> "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Remove the stray empty inline-string cell G5 on sheet 2 ---
$ws2.Range("G5").ClearContents()
